# DesignFirst / Main.xlsx — SAVE commit
# The only data change in this revision is cell C10 on the "Rules" sheet,
# which goes from 18 to 100 (its existing style/format, s="20", is kept).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 100
